$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (monthly -> weekly, first week) ---
$ws.Cells.Item(2,2).Value = 44381
$ws.Cells.Item(2,3).Value = 540.2418001318053

# --- Add weekly rows 3-11, copying row 2's cell formatting first ---
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$rows = @(
  @{r=3;  a=1; b=44388; c=632.3650882947178},
  @{r=4;  a=2; b=44395; c=783.053048224588},
  @{r=5;  a=3; b=44402; c=942.5301754650492},
  @{r=6;  a=4; b=44409; c=1051.262986828038},
  @{r=7;  a=5; b=44416; c=1083.344825274121},
  @{r=8;  a=6; b=44423; c=1051.343280169708},
  @{r=9;  a=7; b=44430; c=975.2416540263248},
  @{r=10; a=8; b=44437; c=859.0627157860984},
  @{r=11; a=9; b=44444; c=704.1046460870393}
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Cells.Item($r,1).Value = $row.a
  $ws.Cells.Item($r,2).Value = $row.b
  $ws.Cells.Item($r,3).Value = $row.c
  $ws.Cells.Item($r,4).Value = "VETINA-Antibiotics -Cat_2"
  $ws.Cells.Item($r,5).Value = "2021-06-27 00:00:00"
  $ws.Cells.Item($r,6).Value = 44381
  $ws.Cells.Item($r,7).Value = 44447
}
